$d = $word.ActiveDocument

# The paragraph we need to touch is the final paragraph of the document body,
# which originally reads: "Git add . (sube todos los cambios)"
# We locate it defensively: prefer the last paragraph, but if its text does not
# contain the expected marker, search for the paragraph that does.
$target = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($target.Range.Text -notlike "*sube todos los cambios*") {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text -like "*sube todos los cambios*") {
            $target = $cand
            break
        }
    }
}

$r = $target.Range

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Rebuild this paragraph (adding gramStart/gramEnd proofing marks around "sube",
# and moving the closing parenthesis into the same run/paragraph), then append
# three more paragraphs describing "git rm --cached" and "git commit -m", each
# with their own proofing marks, finishing with the bookmark ("_GoBack") now
# sitting just before the closing curly quote of the commit message, and a
# trailing empty paragraph.
$xml = "<w:p $ns w:rsidR=`"006E408D`" w:rsidRPr=`"006E408D`" w:rsidRDefault=`"006E408D`"><w:pPr><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr></w:pPr>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>Git</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>add</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> .</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> (</w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>sube</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> todos los cambios)</w:t></w:r>" +
  "</w:p>" +
  "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr></w:pPr>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>Git</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>rm</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> –</w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>cached</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> “archivo” saca el archivo del cache</w:t></w:r>" +
  "</w:p>" +
  "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr></w:pPr>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>Git</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>commit</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> –m “mensaje descriptivo de lo que comite</w:t></w:r>" +
  "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>”</w:t></w:r>" +
  "</w:p>" +
  "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr></w:pPr></w:p>"

$r.InsertXML($xml)
